# LDLC "suivi smartphones" tracker: a new scrape pass happened, so a new
# timestamp column is inserted right before the "nom" / "url_produit"
# columns (old CP -> CQ, old CQ -> CR), and the new CP column gets a
# header timestamp on row 1 plus, for rows that already carry numeric
# pricing history (2-80), the latest price carried forward from CO.
# Rows 81-206 have no pricing history yet, so the new CP cell there
# simply stays blank, matching the other blank history cells on those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at CP (94th column), shifting the existing "nom"
# (old CP) to CQ and "url_produit" (old CQ) to CR.
$ws.Range("CP1").EntireColumn.Insert()

# New header cell CP1 should look like the other timestamp headers
# (bold, bordered, centered) - copy the style from the previous
# timestamp header (now at CO1) then set its own text.
$ws.Range("CO1").Copy()
$ws.Range("CP1").PasteSpecial(-4122)
$ws.Range("CP1").Value2 = "2026-01-31 22:11:23"

# Rows 2-80 already have a numeric price history running through column
# CO; carry that latest price into the freshly inserted CP column too.
$ws.Range("CP2:CP80").Value2 = $ws.Range("CO2:CO80").Value2

# Rows 81-206 have no price history yet (CO is blank there), so leave
# the new CP cells blank as well - nothing further to do for them.
